# Auto-generated Excel COM-interop script applying the Durandal_Profits.xlsx market-data refresh.
# For each changed row, sets the updated cell values (currentAveragePrice / NQ / HQ,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) and clears cells that the refresh removed entirely.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 734.90247
$ws.Range("J17").Value = 730.8857400000001
$ws.Range("L17").Value = 2192.65722
$ws.Range("N17").Value = -2528.65722
# Row 62
$ws.Range("H62").Value = 2776.5386
$ws.Range("I62").Value = 2720
$ws.Range("J62").Value = 2965
$ws.Range("K62").Value = 2720
$ws.Range("L62").Value = 2965
$ws.Range("M62").Value = -2096
$ws.Range("N62").Value = -4213
# Row 65
$ws.Range("H65").Value = 2776.5386
$ws.Range("I65").Value = 2720
$ws.Range("J65").Value = 2965
$ws.Range("K65").Value = 13600
$ws.Range("L65").Value = 14825
$ws.Range("M65").Value = -10480
$ws.Range("N65").Value = -21065
# Row 98
$ws.Range("H98").Value = 54164.95
$ws.Range("I98").Value = 59544.39
$ws.Range("K98").Value = 59544.39
$ws.Range("M98").Value = -58046.39
# Row 122
$ws.Range("H122").Value = 54164.95
$ws.Range("I122").Value = 59544.39
$ws.Range("K122").Value = 178633.17
$ws.Range("M122").Value = -176183.17
# Row 128
$ws.Range("H128").Value = 37424.75
$ws.Range("J128").Value = 37424.75
$ws.Range("L128").Value = 37424.75
$ws.Range("N128").Value = -47384.75
# Row 137
$ws.Range("H137").Value = 1000.2632
$ws.Range("I137").Value = 780.1818
$ws.Range("J137").Value = 1302.875
$ws.Range("K137").Value = 2340.5454
$ws.Range("L137").Value = 3908.625
$ws.Range("M137").Value = 209.4546
$ws.Range("N137").Value = -9008.625
# Row 138
$ws.Range("H138").Value = 2734.9243
$ws.Range("I138").Value = 1427.3864
$ws.Range("J138").Value = 5350
$ws.Range("K138").Value = 4282.1592
$ws.Range("L138").Value = 16050
$ws.Range("M138").Value = 857.8407999999999
$ws.Range("N138").Value = -26330

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 330918.03
$ws.Range("I32").Value = 2092.756
$ws.Range("K32").Value = 2092.756
$ws.Range("M32").Value = -1805.756
# Row 62
$ws.Range("H62").Value = 424.5
$ws.Range("J62").Value = 424.5
$ws.Range("L62").Value = 424.5
$ws.Range("N62").Value = -1672.5
# Row 65
$ws.Range("H65").Value = 424.5
$ws.Range("J65").Value = 424.5
$ws.Range("L65").Value = 1273.5
$ws.Range("N65").Value = -7513.5
# Row 75
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31748
# Row 78
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -98736

$ws = $wb.Worksheets.Item("BSM")
# Row 36
$ws.Range("H36").Value = 873.25
$ws.Range("I36").Value = 873.25
$ws.Range("K36").Value = 873.25
$ws.Range("M36").Value = -339.25
# Row 54
$ws.Range("H54").Value = 10083
$ws.Range("I54").Value = 10083
$ws.Range("K54").Value = 10083
$ws.Range("M54").Value = -9599
# Row 134
$ws.Range("H134").Value = 11773.091
$ws.Range("I134").Value = 5039.1763
$ws.Range("J134").Value = 34668.4
$ws.Range("K134").Value = 15117.5289
$ws.Range("L134").Value = 104005.2
$ws.Range("M134").Value = -12582.5289
$ws.Range("N134").Value = -109075.2

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2645.4075
$ws.Range("I31").Value = 2305.7273
$ws.Range("K31").Value = 2305.7273
$ws.Range("M31").Value = -2010.7273
# Row 34
$ws.Range("H34").Value = 2645.4075
$ws.Range("I34").Value = 2305.7273
$ws.Range("K34").Value = 2305.7273
$ws.Range("M34").Value = -2103.7273
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 724.75
$ws.Range("I92").Value = 660
$ws.Range("J92").Value = 832.6667
$ws.Range("K92").Value = 1980
$ws.Range("L92").Value = 2498.0001
$ws.Range("M92").Value = -732
$ws.Range("N92").Value = -4994.0001
# Row 131
$ws.Range("H131").Value = 22728578
$ws.Range("J131").Value = 23810820
$ws.Range("L131").Value = 71432460
$ws.Range("N131").Value = -71442540
# Row 139
$ws.Range("H139").Value = 6066
$ws.Range("I139").Value = 6757.5
$ws.Range("J139").Value = 3300
$ws.Range("K139").Value = 20272.5
$ws.Range("L139").Value = 9900
$ws.Range("M139").Value = -15132.5
$ws.Range("N139").Value = -20180

$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
# Row 55
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("M55").ClearContents()
# Row 74
$ws.Range("H74").Value = 27233.334
$ws.Range("J74").Value = 27233.334
$ws.Range("L74").Value = 27233.334
$ws.Range("N74").Value = -29105.334
# Row 77
$ws.Range("H77").Value = 27233.334
$ws.Range("J77").Value = 27233.334
$ws.Range("L77").Value = 81700.00199999999
$ws.Range("N77").Value = -91060.00199999999

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 11459.3
$ws.Range("I46").Value = 1584.7142
$ws.Range("J46").Value = 34500
$ws.Range("K46").Value = 1584.7142
$ws.Range("L46").Value = 34500
$ws.Range("M46").Value = -1396.7142
$ws.Range("N46").Value = -34876
# Row 54
$ws.Range("H54").Value = 7000
$ws.Range("J54").Value = 7000
$ws.Range("L54").Value = 7000
$ws.Range("N54").Value = -8288
# Row 55
$ws.Range("H55").Value = 759.7037
$ws.Range("I55").Value = 799.44446
$ws.Range("J55").Value = 680.2222
$ws.Range("K55").Value = 799.44446
$ws.Range("L55").Value = 680.2222
$ws.Range("M55").Value = -626.44446
$ws.Range("N55").Value = -1026.2222
# Row 132
$ws.Range("H132").Value = 3351.2932
$ws.Range("I132").Value = 4923.1665
$ws.Range("J132").Value = 1350.7273
$ws.Range("K132").Value = 14769.4995
$ws.Range("L132").Value = 4052.1819
$ws.Range("M132").Value = -12239.4995
$ws.Range("N132").Value = -9112.1819
# Row 136
$ws.Range("H136").Value = 4904.5454
$ws.Range("I136").Value = 3368.75
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 10106.25
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -7556.25
$ws.Range("N136").Value = -32100

$ws = $wb.Worksheets.Item("WVR")
# Row 13
$ws.Range("H13").Value = 13333
$ws.Range("I13").Value = 20000
$ws.Range("J13").Value = 9999.5
$ws.Range("K13").Value = 20000
$ws.Range("L13").Value = 9999.5
$ws.Range("M13").Value = -19860
$ws.Range("N13").Value = -10279.5
# Row 41
$ws.Range("H41").Value = 14188.5
$ws.Range("J41").Value = 14188.5
$ws.Range("L41").Value = 14188.5
$ws.Range("N41").Value = -14968.5
# Row 132
$ws.Range("H132").Value = 28498598
$ws.Range("I132").Value = 39064044
$ws.Range("J132").Value = 2491351.2
$ws.Range("K132").Value = 117192132
$ws.Range("L132").Value = 7474053.600000001
$ws.Range("M132").Value = -117189602
$ws.Range("N132").Value = -7479113.600000001
